$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after "...during the weekend."
# (the very end of the document). We delete it here and re-insert it at its
# new location (inside the freshly-typed 14th September paragraph) as part
# of the XML fragment below -- this mirrors how Word relocates the auto
# "last edit" bookmark as new text is typed after it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r = $d.Content
$r.Collapse(0)

$xml = @"
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>13</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> of September 2018</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">Today Marcel asked us about our duo app. I’m working together with </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Dylano</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> and our plan is to make an app based on your location. The app will tell you when you discover a new, interesting location. Earlier we asked Chris and Evert their opinions and they both liked it. Evert mentioned that we should focus on what makes us different from other apps. Marcel told us that it was a good idea to start working on our app right away, but we both decided that it would be a better choice for us to do all the smaller assignments first. After the feedback I finished my SMS-History app.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>14</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:vertAlign w:val="superscript"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>th</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> of September 2018</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">Today I worked on the assignments for web service. </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Dylano</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> and me both struggled a bit with this assignment. Both of us managed to make a link between our app and the Firebase database, but only </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Dylano</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> managed to do the Firebase authentication. After following a tutorial my app still doesn’t work. I’m having trouble with the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>onClick</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>-event for some reason and I will ask my teach</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">er about it next week. </w:t>
      </w:r>
    </w:p>
"@

$r.InsertXML($xml)

Write-Host "Paragraphs: $($d.Paragraphs.Count)"
